$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Energy Distribution table (Table2, A6:E13) - Thief column (E), cost-1 row (E9)
$ws.Range("E9").Value = 43

# Update Energy Distribution table (Table24, G7:K12) - Thief column (K), cost-1 & cost-2 rows
$ws.Range("K9").Value = 22
$ws.Range("K10").Value = 34

# Remove the stray N15 helper cell/formula (=75-12) that is no longer needed
$ws.Range("N15").Clear()

# Move the active selection to H2 (matches the saved workbook view state)
$ws.Range("H2").Select()

$wb.Save()
